$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.626.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.499.41"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.560"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.495.07"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0987"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.88%  "

$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.920.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.459.79"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.24"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.485.96"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.74"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.412"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.04%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.37%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.42"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0753"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.57"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.18%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.70"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.33%  "

$ws.Range("E33").Value = "  +2.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.34"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.03"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.68"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.803"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.21"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.47"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.599"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.60"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0912"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0493"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0215"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.19"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.741.80"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.69%  "

